$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21/22 swap: Avalanche moves to row 21, Dai moves to row 22
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"

# Update Price (D) and Volume(1h) (E) columns for rows 2-51
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.197.64"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.52"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.68"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6061"
$ws.Range("E6").Value = "  -3.56%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07124"
$ws.Range("E8").Value = "  -4.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2833"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.09"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07661"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.37"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.791"
$ws.Range("E13").Value = "  -3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6428"
$ws.Range("E14").Value = "  -5.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001009"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.076.63"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.65"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.020"
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.224.93"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.90"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.80"
$ws.Range("E21").Value = "  -4.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.053"
$ws.Range("E23").Value = "  -5.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9971"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.44"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.056"
$ws.Range("E26").Value = "  -4.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1288"
$ws.Range("E27").Value = "  -4.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.73"
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06946"
$ws.Range("E29").Value = "  +6.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.462"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.452"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.837"
$ws.Range("E32").Value = "  -5.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.813"
$ws.Range("E33").Value = "  -6.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.141"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.727"
$ws.Range("E35").Value = "  -6.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6641"
$ws.Range("E36").Value = "  -4.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.532"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.236.72"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.755"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01771"
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.592"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9316"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.999.78"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.24"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.56"
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.645"
$ws.Range("E48").Value = "  -3.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.589"
$ws.Range("E49").Value = "  -6.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05594"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.513"
$ws.Range("E51").Value = "  -5.43%  "
